$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = "[0.07140016525138648, 15.052065762932145]"
$ws.Range("N2").Value = 0.04794514678315087
$ws.Range("O2").Value = 0.04794514678315087
$ws.Range("Q2").Value = "[-2.8680005005744658, -0.07547369738353815]"
$ws.Range("R2").Value = 0.03929058785902484
$ws.Range("S2").Value = 0.03929058785902484
$ws.Range("U2").Value = "[6.065980989041525, 14.764498429972885]"
$ws.Range("V2").Value = [double]"1.653330283701848e-05"
$ws.Range("W2").Value = [double]"1.653330283701848e-05"
$ws.Range("Y2").Value = 0.2642642642642627
$ws.Range("Z2").Value = 10.04204204204204
$ws.Range("M3").Value = "[-0.5555803764303278, 14.835617282299726]"
$ws.Range("N3").Value = 0.06818287571839243
$ws.Range("O3").Value = 0.06818287571839243
$ws.Range("Q3").Value = "[-2.9308952483940813, 0.9685791164220787]"
$ws.Range("R3").Value = 0.3162169293762762
$ws.Range("S3").Value = 0.3162169293762762
$ws.Range("U3").Value = "[5.490561805274887, 13.417136128023]"
$ws.Range("V3").Value = [double]"1.759464668071864e-05"
$ws.Range("W3").Value = [double]"1.759464668071864e-05"
$ws.Range("Y3").Value = -3.391391391391394
$ws.Range("Z3").Value = 10.26226226226226
$ws.Range("M4").Value = "[-1.963656203563417, 12.808929775964067]"
$ws.Range("N4").Value = 0.1461998308724193
$ws.Range("O4").Value = 0.1461998308724193
$ws.Range("Q4").Value = "[-3.4403427057329665, 2.7862373284089657]"
$ws.Range("R4").Value = 0.8333880597101497
$ws.Range("S4").Value = 0.8333880597101497
$ws.Range("U4").Value = "[4.350274134295908, 12.429393079738649]"
$ws.Range("V4").Value = 0.0001313152122215389
$ws.Range("W4").Value = 0.0001313152122215389
$ws.Range("Y4").Value = -9.755755755755757
$ws.Range("Z4").Value = 12.04604604604604
$ws.Range("M5").Value = "[0.14703303445842408, 12.97091041204949]"
$ws.Range("N5").Value = 0.04518099207009008
$ws.Range("O5").Value = 0.04518099207009008
$ws.Range("Q5").Value = "[-1.402552876377424, 1.2767633807381937]"
$ws.Range("R5").Value = 0.9250850436877471
$ws.Range("S5").Value = 0.9250850436877471
$ws.Range("U5").Value = "[4.564487419696486, 12.241969101768959]"
$ws.Range("V5").Value = [double]"6.386250470247212e-05"
$ws.Range("W5").Value = [double]"6.386250470247212e-05"
$ws.Range("Y5").Value = -4.7427627627628
$ws.Range("Z5").Value = 5.210030030030075
$ws.Range("M6").Value = "[-1.0709606776319678, 13.809289267115895]"
$ws.Range("N6").Value = 0.09153864217874164
$ws.Range("O6").Value = 0.09153864217874164
$ws.Range("Q6").Value = "[-1.7296055650394253, 1.478026573760963]"
$ws.Range("R6").Value = 0.875188055739196
$ws.Range("S6").Value = 0.875188055739196
$ws.Range("U6").Value = "[5.144916165928965, 13.59037715189385]"
$ws.Range("V6").Value = [double]"5.278328841051838e-05"
$ws.Range("W6").Value = [double]"5.278328841051838e-05"
$ws.Range("Y6").Value = -5.490390390390439
$ws.Range("Z6").Value = 6.424924924924981
$ws.Range("M7").Value = "[-1.0231935188236676, 13.858208095026786]"
$ws.Range("N7").Value = 0.08920398022203502
$ws.Range("O7").Value = 0.08920398022203502
$ws.Range("Q7").Value = "[-1.459158149415079, 1.7233160902574634]"
$ws.Range("R7").Value = 0.8679785747975655
$ws.Range("S7").Value = 0.8679785747975655
$ws.Range("U7").Value = "[5.1716517160686735, 13.646467106024616]"
$ws.Range("V7").Value = [double]"5.207011133645523e-05"
$ws.Range("W7").Value = [double]"5.207011133645523e-05"
$ws.Range("Y7").Value = 16.9384384384386
$ws.Range("Z7").Value = 28.76030030030056
$ws.Range("M8").Value = "[0.10086682107480094, 12.874711158917389]"
$ws.Range("N8").Value = 0.04663585877134357
$ws.Range("O8").Value = 0.04663585877134357
$ws.Range("Q8").Value = "[-1.2578949563923096, 1.4214213007233099]"
$ws.Range("R8").Value = 0.9027126860165176
$ws.Range("S8").Value = 0.9027126860165176
$ws.Range("U8").Value = "[4.532686903243214, 12.237446043411722]"
$ws.Range("V8").Value = [double]"6.923251313750711e-05"
$ws.Range("W8").Value = [double]"6.923251313750711e-05"
$ws.Range("Y8").Value = 18.05987987988004
$ws.Range("Z8").Value = 28.01267267267292
$ws.Range("M9").Value = "[-0.4298086241198238, 13.494901869885478]"
$ws.Range("N9").Value = 0.06524146373408124
$ws.Range("O9").Value = 0.06524146373408124
$ws.Range("Q9").Value = "[-1.2830528555201557, 1.9245792832802326]"
$ws.Range("R9").Value = 0.6889875392040854
$ws.Range("S9").Value = 0.6889875392040854
$ws.Range("U9").Value = "[4.543467351674554, 12.201739539743127]"
$ws.Range("V9").Value = [double]"6.490877960896313e-05"
$ws.Range("W9").Value = [double]"6.490877960896313e-05"
$ws.Range("Y9").Value = 16.19081081081096
$ws.Range("Z9").Value = 28.10612612612638
$ws.Range("M10").Value = "[-1.0466426361875136, 14.044604119557945]"
$ws.Range("N10").Value = 0.08963421080590517
$ws.Range("O10").Value = 0.08963421080590517
$ws.Range("Q10").Value = "[-1.3962634015954638, 1.8868424345884645]"
$ws.Range("R10").Value = 0.7648327665665189
$ws.Range("S10").Value = 0.7648327665665189
$ws.Range("U10").Value = "[5.191388062648699, 13.5888891102137]"
$ws.Range("V10").Value = [double]"4.69316930824526e-05"
$ws.Range("W10").Value = [double]"4.69316930824526e-05"
$ws.Range("Y10").Value = 16.33099099099114
$ws.Range("Z10").Value = 28.52666666666692
